$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44291
$ws.Range("J2").Value = 30

$ws.Range("D3").Value = 44284
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 10000
$ws.Range("P3").Value = 500

$ws.Range("D4").Value = 44277
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 11000
$ws.Range("L4").Value = 11000
$ws.Range("M4").Value = 11000
$ws.Range("P4").Value = 550

$ws.Range("D5").Value = 44280
